# Apply the per-cell odds/value corrections described in the commit diff.
# Each Range(...).Value assignment below updates a single cell to its new
# value; rows are processed in ascending order, left-to-right by column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("L3").Value = 1.2
$ws.Range("M3").Value = 4.5

# Row 5
$ws.Range("K5").Value = 19

# Row 7
$ws.Range("K7").Value = 13

# Row 9
$ws.Range("G9").Value = 2
$ws.Range("I9").Value = 3.6
$ws.Range("L9").Value = 1.25
$ws.Range("M9").Value = 3.75
$ws.Range("N9").Value = 1.9
$ws.Range("O9").Value = 1.95
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 8
$ws.Range("V9").Value = 9
$ws.Range("AD9").Value = 201
$ws.Range("AG9").Value = 12

# Row 12
$ws.Range("L12").Value = 1.2
$ws.Range("M12").Value = 4.33
$ws.Range("N12").Value = 1.7
$ws.Range("O12").Value = 2.1

# Row 13
$ws.Range("G13").Value = 2.88
$ws.Range("I13").Value = 2.5
$ws.Range("J13").Value = 1.06
$ws.Range("K13").Value = 10
$ws.Range("N13").Value = 1.93
$ws.Range("O13").Value = 1.93
$ws.Range("V13").Value = 11
$ws.Range("Z13").Value = 10
$ws.Range("AE13").Value = 9
$ws.Range("AH13").Value = 23
$ws.Range("AI13").Value = 19

# Row 14 (G,H,I change; J,K become blank; L..AJ all updated)
$ws.Range("G14").Value = 1.25
$ws.Range("H14").Value = 7.5
$ws.Range("I14").Value = 7.5
$ws.Range("J14").ClearContents()
$ws.Range("K14").ClearContents()
$ws.Range("L14").Value = 1.03
$ws.Range("M14").Value = 17
$ws.Range("N14").Value = 1.11
$ws.Range("O14").Value = 6.5
$ws.Range("P14").Value = 1.1
$ws.Range("Q14").Value = 7
$ws.Range("R14").Value = 1.33
$ws.Range("S14").Value = 3.25
$ws.Range("T14").Value = 26
$ws.Range("U14").Value = 17
$ws.Range("V14").Value = 12
$ws.Range("W14").Value = 15
$ws.Range("X14").Value = 11
$ws.Range("Y14").Value = 15
$ws.Range("Z14").Value = 34
$ws.Range("AA14").Value = 21
$ws.Range("AB14").Value = 19
$ws.Range("AC14").Value = 29
$ws.Range("AD14").Value = 67
$ws.Range("AE14").Value = 51
$ws.Range("AF14").Value = 51
$ws.Range("AG14").Value = 26
$ws.Range("AH14").Value = 101
$ws.Range("AI14").Value = 41
$ws.Range("AJ14").Value = 34

# Row 17
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 10

# Row 18
$ws.Range("J18").Value = 1.11
$ws.Range("K18").Value = 6.5
$ws.Range("Z18").Value = 6.5

# Row 19
$ws.Range("G19").Value = 1.9
$ws.Range("I19").Value = 3.9
$ws.Range("K19").Value = 12
$ws.Range("P19").Value = 1.36
$ws.Range("Q19").Value = 3
$ws.Range("Z19").Value = 12

# Row 22
$ws.Range("G22").Value = 1.85
$ws.Range("H22").Value = 3.6
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 1.07
$ws.Range("K22").Value = 9
$ws.Range("U22").Value = 8.5
$ws.Range("X22").Value = 17
$ws.Range("Z22").Value = 9

# Row 25
$ws.Range("K25").Value = 9
$ws.Range("N25").Value = 2.1
$ws.Range("O25").Value = 1.7

# Row 29
$ws.Range("J29").Value = 1.02
$ws.Range("L29").Value = 1.19
